$d = $word.ActiveDocument

# Helper: force a run split at a given character offset by briefly adding
# and then deleting a bookmark there. The OOXML serializer merges adjacent
# runs that share formatting, but a bookmark boundary (even transient)
# keeps the surrounding text in separate <w:r> elements once split.
function Split-RunAt($offset) {
    $sr = $d.Range($offset, $offset)
    $d.Bookmarks.Add("TEMP_SPLIT_MARK", $sr) | Out-Null
    $d.Bookmarks("TEMP_SPLIT_MARK").Delete()
}

# ---------------------------------------------------------------------
# "May:" section - remove the old (hidden) _GoBack bookmark that sits
# inside "Detect two di|screte levels" before editing that text, and
# update "levels" -> "discrete intensities" while keeping the leading
# "-" run and trailing " of knocks..." run intact.
# ---------------------------------------------------------------------
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

$rPrefix = $d.Content
$null = $rPrefix.Find.Execute("Detect two di")
$detectStart = $rPrefix.Start

$rLevels = $d.Content
$null = $rLevels.Find.Execute("screte levels")
$rLevels.Text = "screte intensities"

Split-RunAt ($detectStart + "Detect two discrete ".Length)
Split-RunAt ($detectStart + "Detect two discrete intensities".Length)

# ---------------------------------------------------------------------
# "Must:" section - "-Be capable of visual and audio ques." becomes two
# runs, and the _GoBack bookmark now lands at the end of this paragraph.
# ---------------------------------------------------------------------
$rQues = $d.Content
$null = $rQues.Find.Execute("-Be capable of visual and audio ques.")
$quesStart = $rQues.Start
$rQues.Text = "-Be capable of visual or audio queues for confirmation."
$quesTextEnd = $quesStart + "-Be capable of visual or audio queues for confirmation.".Length

Split-RunAt ($quesStart + "-Be capable of visual ".Length)

# Place the bookmark right after the final period, before the paragraph
# mark. A zero-length range built exactly at story-end-of-text has a
# positioning quirk, so pad with a placeholder character, bookmark the
# gap, then remove the placeholder again.
$phRange = $d.Range($quesTextEnd, $quesTextEnd)
$phRange.InsertAfter("X")
$bmRange = $d.Range($quesTextEnd, $quesTextEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$d.Range($quesTextEnd, $quesTextEnd + 1).Text = ""

# ---------------------------------------------------------------------
# "Should:" section - "-Require user to confirm knock sequence prior to
# writing to EEPROM" becomes "-Require user to confirm knock sequence
# prior to saving", split into two runs.
# ---------------------------------------------------------------------
$rReq = $d.Content
$null = $rReq.Find.Execute("-Require user to confirm knock sequence prior to writing to EEPROM")
$reqStart = $rReq.Start
$rReq.Text = "-Require user to confirm knock sequence prior to saving"

Split-RunAt ($reqStart + "-Require user to ".Length)

# ---------------------------------------------------------------------
# "May:" section - "-Have battery-back up capabilities." becomes three
# runs, then two new empty paragraphs are appended at the end.
# ---------------------------------------------------------------------
$rBatt = $d.Content
$null = $rBatt.Find.Execute("-Have battery-back up capabilities.")
$battStart = $rBatt.Start
$rBatt.Text = "-Have battery-back up and recharging capabilities."

Split-RunAt ($battStart + "-Have battery-back up ".Length)
Split-RunAt ($battStart + "-Have battery-back up and recharging ".Length)

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.Text = "`r"
$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.Text = "`r"
